$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 4 (Ahmed Gamal -> Basem hamdy) ---
# Order matters: it controls the append order of new shared strings so
# the saved sharedStrings.xml table matches the expected layout.
$ws.Range("A4").Value = "Basem"
$ws.Range("C4").Value = "basemhamdy1983@gmail.com"
$ws.Range("B4").Value = "hamdy"
$ws.Range("D4").Value = "P@ssw0rd85"

# --- New columns: Job Title (G), Company Name (H), Work Industry (I) ---
$ws.Range("G1").Value = "Job Title"
$ws.Range("H1").Value = "Company Name"
$ws.Range("G2").Value = "QC Engineer"
$ws.Range("H2").Value = "Arcom"
$ws.Range("G3").Value = "Development Team Lead"
$ws.Range("H3").Value = "Arcom"
$ws.Range("G4").Value = "Manager"
$ws.Range("H4").Value = "Arcom"
$ws.Range("I1").Value = "Work Industry"
$ws.Range("I4").Value = "Information Technology and Services"
$ws.Range("I3").Value = "Information Services"
$ws.Range("I2").Value = "Hospital & Health Care"

# New cell style for I4: small Consolas font
$ws.Range("I4").Font.Name = "Consolas"
$ws.Range("I4").Font.Size = 9
$ws.Range("I4").Font.Color = 2236962

# --- Hyperlinks: drop the (now plain) C4 e-mail link, keep the rest ---
# Deleting via the Hyperlinks collection clears the whole sheet, so
# rebuild the ones that should survive (C2, D2, C3, D3, D4).
$ws.Range("C4").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:islamtalkha83@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:P@ssw0rd85")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:ali89afit@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:P@ssw0rd85")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:P@ssw0rd85")

# --- Column sizing for the new / widened columns ---
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(9).AutoFit()

# --- Page setup / selection ---
$ws.PageSetup.Orientation = 1
$ws.Range("K3").Select()
